$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Objetivos:" row (row 10) — replace the long objectives text with the
# teacher identifier that now belongs there.
$ws.Cells.Item(10, 2).Value = "8822123 - Roberta Veloso Garcia"
$ws.Cells.Item(10, 3).Value = "8822123 - Roberta Veloso Garcia"

# Row 13 held only the "8822123 - Roberta Veloso Garcia" value in B/C with
# an empty A cell (no label) — it is removed entirely, shifting everything
# below it up by one row.
$ws.Rows.Item(13).Delete()

# After the shift, "Programa resumido:" is now row 13 — its value becomes
# "Semestral".
$ws.Cells.Item(13, 2).Value = "Semestral"
$ws.Cells.Item(13, 3).Value = "Semestral"

# "Programa:" is now row 15 — its value becomes "01/01/2018". Copy the
# existing "01/01/2018" text cell (row 8) and paste-special values-only so
# the string stays plain text instead of being auto-parsed into a date
# serial number (which would change the cell's type/style).
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4163)
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4163)

# "Método:" is now row 18 — its value becomes the teacher identifier.
$ws.Cells.Item(18, 2).Value = "8822123 - Roberta Veloso Garcia"
$ws.Cells.Item(18, 3).Value = "8822123 - Roberta Veloso Garcia"

# "Critério:" is now row 19 — its value becomes the "NF=A avaliação..." text.
$ws.Cells.Item(19, 2).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Cells.Item(19, 3).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# "Norma de recuperação:" is now row 20 — its value becomes "NF≥ 5,0.".
$ws.Cells.Item(20, 2).Value = "NF≥ 5,0."
$ws.Cells.Item(20, 3).Value = "NF≥ 5,0."

# "Bibliografia:" is now row 21 — its value becomes the "(NF+RC)/2..." text.
$ws.Cells.Item(21, 2).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Cells.Item(21, 3).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
